$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Ref, $Val)
    $cell = $ws.Range($Ref)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = $origStyle
}

Set-TextValue "D2" "28.854.46"
Set-TextValue "E2" "  +7.87%  "
Set-TextValue "D3" "1.816.23"
Set-TextValue "E3" "  +5.31%  "
Set-TextValue "D4" "0.9994"
Set-TextValue "E4" "  +0.15%  "
Set-TextValue "D5" "251.19"
Set-TextValue "E5" "  +4.00%  "
Set-TextValue "D6" "0.9994"
Set-TextValue "E6" "  +0.08%  "
Set-TextValue "D7" "0.4977"
Set-TextValue "E7" "  +2.30%  "
Set-TextValue "D8" "0.2802"
Set-TextValue "E8" "  +8.50%  "
Set-TextValue "D9" "0.06423"
Set-TextValue "E9" "  +3.61%  "
Set-TextValue "D10" "1.818.71"
Set-TextValue "E10" "  +5.45%  "
Set-TextValue "D11" "16.83"
Set-TextValue "E11" "  +5.52%  "
Set-TextValue "D12" "0.07174"
Set-TextValue "E12" "  +3.70%  "
Set-TextValue "D13" "0.6533"
Set-TextValue "E13" "  +7.48%  "
Set-TextValue "D14" "84.23"
Set-TextValue "E14" "  +9.66%  "
Set-TextValue "D15" "4.739"
Set-TextValue "E15" "  +5.87%  "
Set-TextValue "D16" "28.843.23"
Set-TextValue "E16" "  +8.60%  "
Set-TextValue "D17" "0.9985"
Set-TextValue "E17" "  -0.02%  "
Set-TextValue "D18" "0.000007437"
Set-TextValue "E18" "  +3.71%  "
Set-TextValue "D19" "0.9995"
Set-TextValue "E19" "  +0.14%  "
Set-TextValue "D20" "12.33"
Set-TextValue "E20" "  +7.91%  "
Set-TextValue "D21" "2.051.79"
Set-TextValue "E21" "  +5.05%  "
Set-TextValue "D22" "4.632"
Set-TextValue "E22" "  +4.75%  "
Set-TextValue "D23" "8.931"
Set-TextValue "E23" "  +4.35%  "
Set-TextValue "D24" "5.374"
Set-TextValue "E24" "  +6.16%  "
Set-TextValue "D25" "143.98"
Set-TextValue "E25" "  +4.98%  "
Set-TextValue "D26" "132.81"
Set-TextValue "E26" "  +25.08%  "
Set-TextValue "D27" "16.46"
Set-TextValue "E27" "  +7.83%  "
Set-TextValue "E28" "  +7.45%  "
Set-TextValue "D29" "1.407"
Set-TextValue "E29" "  +1.39%  "
Set-TextValue "D30" "4.183"
Set-TextValue "E30" "  +6.42%  "
Set-TextValue "D31" "0.08401"
Set-TextValue "E31" "  +5.67%  "
Set-TextValue "D32" "3.887"
Set-TextValue "E32" "  +5.40%  "
Set-TextValue "D33" "0.04979"
Set-TextValue "E33" "  +10.89%  "
Set-TextValue "D34" "1.094"
Set-TextValue "E34" "  +8.43%  "
Set-TextValue "B35" "ImmutableX"
Set-TextValue "C35" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D35" "0.6841"
Set-TextValue "E35" "  +9.90%  "
Set-TextValue "B36" "HuobiToken"
Set-TextValue "C36" "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "2.695"
Set-TextValue "E36" "  +3.79%  "
Set-TextValue "B37" "MXToken"
Set-TextValue "C37" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D37" "2.790"
Set-TextValue "E37" "  +14.49%  "
Set-TextValue "B38" "RenderToken"
Set-TextValue "C38" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D38" "2.242"
Set-TextValue "E38" "  +9.98%  "
Set-TextValue "B39" "TrustWalletToken"
Set-TextValue "C39" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D39" "0.9698"
Set-TextValue "E39" "  +4.64%  "
Set-TextValue "B40" "FraxShare"
Set-TextValue "C40" "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D40" "6.086"
Set-TextValue "E40" "  +7.91%  "
Set-TextValue "D41" "0.01597"
Set-TextValue "E41" "  +6.76%  "
Set-TextValue "B42" "PaxDollar"
Set-TextValue "C42" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D42" "0.9994"
Set-TextValue "E42" "  +0.13%  "
Set-TextValue "B43" "TheSandbox"
Set-TextValue "C43" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D43" "0.4132"
Set-TextValue "E43" "  +7.86%  "
Set-TextValue "D44" "101.10"
Set-TextValue "E44" "  +1.61%  "
Set-TextValue "B45" "Aptos"
Set-TextValue "C45" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D45" "7.278"
Set-TextValue "E45" "  +6.33%  "
Set-TextValue "B46" "Algorand"
Set-TextValue "C46" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D46" "0.1230"
Set-TextValue "E46" "  +6.26%  "
Set-TextValue "B47" "Cronos"
Set-TextValue "C47" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D47" "0.05525"
Set-TextValue "E47" "  +2.43%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "8.205"
Set-TextValue "E48" "  +4.59%  "
Set-TextValue "B49" "Elrond"
Set-TextValue "C49" "https://coinranking.com/coin/omwkOTglq+elrond-egld"
Set-TextValue "D49" "31.82"
Set-TextValue "E49" "  +5.66%  "
Set-TextValue "B50" "Decentraland"
Set-TextValue "C50" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D50" "0.3657"
Set-TextValue "E50" "  +8.68%  "
Set-TextValue "B51" "NEARProtocol"
Set-TextValue "C51" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D51" "1.315"
Set-TextValue "E51" "  +7.03%  "
